{"js": "// Apply the text edits described by the diff:\n//  - Update the H1 title (and its later bold restatement) from\n//    \"Play Joker Poker MH Free - Simple and Clean Design\" to\n//    \"Play Joker Poker MH for Free\"\n//  - Rewrite the \"What we like\" bullet list items\n//  - Rewrite one \"What we don't like\" bullet list item\n//  - Rewrite the closing italic summary line\n//\n// Uses Body.search + Range.insertText(\u2026, \"Replace\") so that existing\n// (possibly empty) sibling runs in the paragraph are left alone and only\n// the matched text range is swapped out.\n\nconst replacements = [\n  [\n    \"Play Joker Poker MH Free - Simple and Clean Design\",\n    \"Play Joker Poker MH for Free\",\n  ],\n  [\n    \"Simple and essential gameplay\",\n    \"Simple gameplay and clean design\",\n  ],\n  [\n    \"Clean and user-friendly design\",\n    \"Joker card increases chances of winning combinations\",\n  ],\n  [\n    \"Bonus game that multiplies winnings\",\n    \"Bonus game adds excitement and potential for higher winnings\",\n  ],\n  [\n    \"Demo version available to refine skills\",\n    \"Demo version available for practice before betting real money\",\n  ],\n  [\n    \"No progressive jackpot to win\",\n    \"Limited betting range (1 to 5 coins)\",\n  ],\n  [\n    \"Read our review of Joker Poker MH, a simple online slot game with clean design. Play for free and multiply your winnings with the bonus game.\",\n    \"Read our review of Joker Poker MH and play for free to hone your skills.\",\n  ],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text edits described by the diff:\n#  - Update the H1 title (and its later bold restatement) from\n#    \"Play Joker Poker MH Free - Simple and Clean Design\" to\n#    \"Play Joker Poker MH for Free\"\n#  - Rewrite the \"What we like\" bullet list items\n#  - Rewrite one \"What we don't like\" bullet list item\n#  - Rewrite the closing italic summary line\n#\n# Plain `Range.Text = ...` / `Range.Delete()` on this host normalizes a\n# paragraph's run list (it drops the zero-length `<w:r/>` runs that sit next\n# to the real text run here). Wrapping the matched Range in a temporary\n# bookmark before mutating its `.Text`, then deleting the bookmark (which\n# only removes the bookmark markers, not the content), keeps the edit\n# scoped to just that run's text and leaves sibling runs untouched.\n\nfunction Replace-AllOccurrences {\n    param(\n        $doc,\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $bookmarkName = \"zzTmpReplace\"\n    $replacedCount = 0\n\n    while ($true) {\n        $rng = $doc.Content\n        $find = $rng.Find\n        $find.ClearFormatting()\n        $find.Text = $FindText\n        $find.MatchCase = $true\n        $find.MatchWholeWord = $false\n        $find.MatchWildcards = $false\n        $find.Forward = $true\n        $find.Wrap = 0  # wdFindStop: do not wrap, so each loop iteration finds the next fresh match\n\n        $found = $find.Execute()\n        if (-not $found) {\n            break\n        }\n\n        # $rng has been mutated by Find.Execute to cover the matched text.\n        if ($doc.Bookmarks.Exists($bookmarkName)) {\n            $doc.Bookmarks($bookmarkName).Delete()\n        }\n        $doc.Bookmarks.Add($bookmarkName, $rng)\n        $bmRange = $doc.Bookmarks($bookmarkName).Range\n        $bmRange.Text = $ReplaceText\n        $doc.Bookmarks($bookmarkName).Delete()\n\n        $replacedCount = $replacedCount + 1\n        if ($replacedCount -gt 50) {\n            break\n        }\n    }\n\n    return $replacedCount\n}\n\n$d = $word.ActiveDocument\n\nReplace-AllOccurrences $d \"Play Joker Poker MH Free - Simple and Clean Design\" \"Play Joker Poker MH for Free\" | Out-Null\nReplace-AllOccurrences $d \"Simple and essential gameplay\" \"Simple gameplay and clean design\" | Out-Null\nReplace-AllOccurrences $d \"Clean and user-friendly design\" \"Joker card increases chances of winning combinations\" | Out-Null\nReplace-AllOccurrences $d \"Bonus game that multiplies winnings\" \"Bonus game adds excitement and potential for higher winnings\" | Out-Null\nReplace-AllOccurrences $d \"Demo version available to refine skills\" \"Demo version available for practice before betting real money\" | Out-Null\nReplace-AllOccurrences $d \"No progressive jackpot to win\" \"Limited betting range (1 to 5 coins)\" | Out-Null\nReplace-AllOccurrences $d \"Read our review of Joker Poker MH, a simple online slot game with clean design. Play for free and multiply your winnings with the bonus game.\" \"Read our review of Joker Poker MH and play for free to hone your skills.\" | Out-Null\n"}
